$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 675.1429000000001
$ws.Range("I19").Value = 410
$ws.Range("J19").Value = 747.4545000000001
$ws.Range("K19").Value = 410
$ws.Range("L19").Value = 747.4545000000001
$ws.Range("M19").Value = -235
$ws.Range("N19").Value = -1097.4545
# row 62
$ws.Range("H62").Value = 3095843.5
$ws.Range("I62").Value = 4837581
$ws.Range("J62").Value = 14307.692
$ws.Range("K62").Value = 4837581
$ws.Range("L62").Value = 14307.692
$ws.Range("M62").Value = -4836957
$ws.Range("N62").Value = -15555.692
# row 65
$ws.Range("H65").Value = 3095843.5
$ws.Range("I65").Value = 4837581
$ws.Range("J65").Value = 14307.692
$ws.Range("K65").Value = 24187905
$ws.Range("L65").Value = 71538.45999999999
$ws.Range("M65").Value = -24184785
$ws.Range("N65").Value = -77778.45999999999
# row 106
$ws.Range("H106").Value = 5053710.5
$ws.Range("I106").Value = 5851038.5
$ws.Range("J106").Value = 3966.6667
$ws.Range("K106").Value = 5851038.5
$ws.Range("L106").Value = 3966.6667
$ws.Range("M106").Value = -5850407.5
$ws.Range("N106").Value = -5228.6667
# row 107
$ws.Range("H107").Value = 505420.62
$ws.Range("I107").Value = 505420.62
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 505420.62
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -503500.62
# row 132
$ws.Range("H132").Value = 275520.72
$ws.Range("I132").Value = 290414.9
$ws.Range("K132").Value = 871244.7000000001
$ws.Range("M132").Value = -868714.7000000001
# row 137
$ws.Range("H137").Value = 23256922
$ws.Range("I137").Value = 30303824
$ws.Range("J137").Value = 2147.8
$ws.Range("K137").Value = 90911472
$ws.Range("L137").Value = 6443.400000000001
$ws.Range("M137").Value = -90908922
$ws.Range("N137").Value = -11543.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 16971.395
$ws.Range("I32").Value = 2654.8413
$ws.Range("K32").Value = 2654.8413
$ws.Range("M32").Value = -2367.8413
# row 61
$ws.Range("H61").Value = 1725.3654
$ws.Range("I61").Value = 1288.3778
$ws.Range("J61").Value = 4534.5713
$ws.Range("K61").Value = 1288.3778
$ws.Range("L61").Value = 4534.5713
$ws.Range("M61").Value = -1076.3778
$ws.Range("N61").Value = -4958.5713
# row 74
$ws.Range("H74").Value = 3103.3582
$ws.Range("I74").Value = 1027.9387
$ws.Range("J74").Value = 8753.111000000001
$ws.Range("K74").Value = 1027.9387
$ws.Range("L74").Value = 8753.111000000001
$ws.Range("M74").Value = -153.9386999999999
$ws.Range("N74").Value = -10501.111
# row 77
$ws.Range("H77").Value = 3103.3582
$ws.Range("I77").Value = 1027.9387
$ws.Range("J77").Value = 8753.111000000001
$ws.Range("K77").Value = 5139.693499999999
$ws.Range("L77").Value = 43765.55500000001
$ws.Range("M77").Value = -771.6934999999994
$ws.Range("N77").Value = -52501.55500000001
# row 122
$ws.Range("H122").Value = 8667.0625
$ws.Range("I122").Value = 9211.532999999999
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 27634.599
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -25184.599
$ws.Range("N122").Value = -6400
# row 132
$ws.Range("H132").Value = 2693.1892
$ws.Range("I132").Value = 2161.8125
$ws.Range("J132").Value = 6094
$ws.Range("K132").Value = 6485.4375
$ws.Range("L132").Value = 18282
$ws.Range("M132").Value = -3955.4375
$ws.Range("N132").Value = -23342
# row 136
$ws.Range("H136").Value = 1725.3654
$ws.Range("I136").Value = 1288.3778
$ws.Range("J136").Value = 4534.5713
$ws.Range("K136").Value = 3865.1334
$ws.Range("L136").Value = 13603.7139
$ws.Range("M136").Value = -1315.1334
$ws.Range("N136").Value = -18703.7139
# row 139
$ws.Range("H139").Value = 46359
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 46359
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = 46359
$ws.Range("N139").Value = -56639

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 1289.3143
$ws.Range("I99").Value = 1238.3125
$ws.Range("J99").Value = 1833.3334
$ws.Range("K99").Value = 1238.3125
$ws.Range("L99").Value = 1833.3334
$ws.Range("M99").Value = 259.6875
$ws.Range("N99").Value = -4829.3334
# row 105
$ws.Range("H105").Value = 3224.4
$ws.Range("I105").Value = 3144.389
$ws.Range("K105").Value = 3144.389
$ws.Range("M105").Value = -1397.389
# row 107
$ws.Range("H107").Value = 1179
$ws.Range("I107").Value = 1130
$ws.Range("J107").Value = 1375
$ws.Range("K107").Value = 1130
$ws.Range("L107").Value = 1375
$ws.Range("M107").Value = 790
$ws.Range("N107").Value = -5215
# row 134
$ws.Range("H134").Value = 20410418
$ws.Range("I134").Value = 23257506
$ws.Range("J134").Value = 6302
$ws.Range("K134").Value = 69772518
$ws.Range("L134").Value = 18906
$ws.Range("M134").Value = -69769983
$ws.Range("N134").Value = -23976

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1161.4286
$ws.Range("I16").Value = 546.6667
$ws.Range("J16").Value = 1622.5
$ws.Range("K16").Value = 546.6667
$ws.Range("L16").Value = 1622.5
$ws.Range("M16").Value = -259.6667
$ws.Range("N16").Value = -2196.5
# row 22
$ws.Range("H22").Value = 579.2308
$ws.Range("J22").Value = 764.8333
$ws.Range("L22").Value = 764.8333
$ws.Range("N22").Value = -1464.8333
# row 31
$ws.Range("H31").Value = 1450.0435
$ws.Range("I31").Value = 857.3036
$ws.Range("J31").Value = 4003.3845
$ws.Range("K31").Value = 857.3036
$ws.Range("L31").Value = 4003.3845
$ws.Range("M31").Value = -562.3036
$ws.Range("N31").Value = -4593.3845
# row 34
$ws.Range("H34").Value = 1450.0435
$ws.Range("I34").Value = 857.3036
$ws.Range("J34").Value = 4003.3845
$ws.Range("K34").Value = 857.3036
$ws.Range("L34").Value = 4003.3845
$ws.Range("M34").Value = -655.3036
$ws.Range("N34").Value = -4407.3845
# row 113
$ws.Range("H113").Value = 1161.4286
$ws.Range("I113").Value = 546.6667
$ws.Range("J113").Value = 1622.5
$ws.Range("K113").Value = 546.6667
$ws.Range("L113").Value = 1622.5
$ws.Range("M113").Value = 1623.3333
$ws.Range("N113").Value = -5962.5
# row 132
$ws.Range("H132").Value = 2131.9434
$ws.Range("I132").Value = 1505.738
$ws.Range("J132").Value = 4522.909
$ws.Range("K132").Value = 4517.214
$ws.Range("L132").Value = 13568.727
$ws.Range("M132").Value = -1987.214
$ws.Range("N132").Value = -18628.727
# row 134
$ws.Range("H134").Value = 1893.1904
$ws.Range("I134").Value = 1285.1608
$ws.Range("J134").Value = 6757.4287
$ws.Range("K134").Value = 3855.4824
$ws.Range("L134").Value = 20272.2861
$ws.Range("M134").Value = -1320.4824
$ws.Range("N134").Value = -25342.2861

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 64
$ws.Range("H64").Value = 4363
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 4843.4287
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 14530.2861
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -15070.2861
# row 67
$ws.Range("H67").Value = 4363
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 4843.4287
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 14530.2861
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -16402.2861
# row 140
$ws.Range("H140").Value = 4132.3555
$ws.Range("I140").Value = 5239.76
$ws.Range("J140").Value = 2748.1
$ws.Range("K140").Value = 15719.28
$ws.Range("L140").Value = 8244.299999999999
$ws.Range("M140").Value = -10539.28
$ws.Range("N140").Value = -18604.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Range("H122").Value = 1588688.1
$ws.Range("I122").Value = 2223563.5
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 6670690.5
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -6668240.5
$ws.Range("N122").Value = -9400
# row 123
$ws.Range("H123").Value = 9268.333000000001
$ws.Range("J123").Value = 9268.333000000001
$ws.Range("L123").Value = 9268.333000000001
$ws.Range("N123").Value = -14168.333
# row 132
$ws.Range("H132").Value = 3023.5098
$ws.Range("I132").Value = 2591.8
$ws.Range("J132").Value = 4593.364
$ws.Range("K132").Value = 7775.400000000001
$ws.Range("L132").Value = 13780.092
$ws.Range("M132").Value = -5245.400000000001
$ws.Range("N132").Value = -18840.092

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 14122.5
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 18496.666
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 18496.666
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -19086.666
# row 27
$ws.Range("H27").Value = 14122.5
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 18496.666
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 18496.666
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -18710.666
# row 29
$ws.Range("H29").Value = 17125
$ws.Range("J29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("N29").Value = -8590
# row 46
$ws.Range("H46").Value = 2468.5715
$ws.Range("I46").Value = 2133.3333
$ws.Range("J46").Value = 2720
$ws.Range("K46").Value = 2133.3333
$ws.Range("L46").Value = 2720
$ws.Range("M46").Value = -1945.3333
$ws.Range("N46").Value = -3096
# row 122
$ws.Range("H122").Value = 3544
$ws.Range("I122").Value = 1902
$ws.Range("K122").Value = 5706
$ws.Range("M122").Value = -3256
# row 132
$ws.Range("H132").Value = 6781.3335
$ws.Range("I132").Value = 7624.609
$ws.Range("J132").Value = 4841.8
$ws.Range("K132").Value = 22873.827
$ws.Range("L132").Value = 14525.4
$ws.Range("M132").Value = -20343.827
$ws.Range("N132").Value = -19585.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 14
$ws.Range("H14").Value = 26875.25
$ws.Range("J14").Value = 25333.6
$ws.Range("L14").Value = 25333.6
$ws.Range("N14").Value = -25669.6
# row 126
$ws.Range("H126").Value = 37793.223
$ws.Range("I126").Value = 50665.85
$ws.Range("J126").Value = 1014.2857
$ws.Range("K126").Value = 151997.55
$ws.Range("L126").Value = 3042.8571
$ws.Range("M126").Value = -149527.55
$ws.Range("N126").Value = -7982.8571
# row 128
$ws.Range("H128").Value = 73367.86
$ws.Range("J128").Value = 73367.86
$ws.Range("L128").Value = 73367.86
$ws.Range("N128").Value = -83327.86
# row 132
$ws.Range("H132").Value = 7464437
$ws.Range("I132").Value = 10205731
$ws.Range("J132").Value = 2025.5
$ws.Range("K132").Value = 30617193
$ws.Range("L132").Value = 6076.5
$ws.Range("M132").Value = -30614663
$ws.Range("N132").Value = -11136.5
# row 136
$ws.Range("H136").Value = 17698.4
$ws.Range("I136").Value = 19091.62
$ws.Range("J136").Value = 2373
$ws.Range("K136").Value = 57274.86
$ws.Range("L136").Value = 7119
$ws.Range("M136").Value = -54724.86
$ws.Range("N136").Value = -12219
